$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label
$ws.Range("F1").Value = "T Tag"

# Map of full speaker names/labels to their abbreviated forms
$map = @{
    "ANTOINETTE VILLARIN" = "T"
    "TEACHER"             = "T"
    "STUDENT"             = "S"
    "STUDENTS"            = "SS"
}

for ($r = 2; $r -le 143; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
